# Update Resumo de Inscricoes values per the latest processing run.
# Commit via gitrun.py em 2024-09-18 12:00:42
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 59
$ws.Range("H2").Value = 59
$ws.Range("F3").Value = 26
$ws.Range("H3").Value = 26
$ws.Range("F5").Value = 77
$ws.Range("H5").Value = 77
$ws.Range("E10").Value = 403
$ws.Range("F10").Value = 196
$ws.Range("H10").Value = 196
$ws.Range("E11").Value = 273
$ws.Range("F11").Value = 147
$ws.Range("H11").Value = 147
$ws.Range("E12").Value = 388
$ws.Range("F12").Value = 214
$ws.Range("H12").Value = 214
$ws.Range("E13").Value = 103
$ws.Range("F13").Value = 56
$ws.Range("H13").Value = 56
$ws.Range("F14").Value = 50
$ws.Range("H14").Value = 50
$ws.Range("E15").Value = 128
$ws.Range("F15").Value = 51
$ws.Range("H15").Value = 51
$ws.Range("E16").Value = 169
$ws.Range("F16").Value = 85
$ws.Range("H16").Value = 85
$ws.Range("F17").Value = 36
$ws.Range("H17").Value = 36
$ws.Range("E20").Value = 77
$ws.Range("F21").Value = 68
$ws.Range("H21").Value = 68
$ws.Range("F22").Value = 72
$ws.Range("H22").Value = 72
$ws.Range("F23").Value = 74
$ws.Range("H23").Value = 74
$ws.Range("E24").Value = 172
$ws.Range("F24").Value = 92
$ws.Range("H24").Value = 92
$ws.Range("E25").Value = 201
$ws.Range("F25").Value = 93
$ws.Range("H25").Value = 93
$ws.Range("E26").Value = 116
$ws.Range("F26").Value = 71
$ws.Range("H26").Value = 71
$ws.Range("F27").Value = 127
$ws.Range("H27").Value = 127
$ws.Range("E28").Value = 156
$ws.Range("F28").Value = 53
$ws.Range("H28").Value = 53
$ws.Range("E30").Value = 172
$ws.Range("F30").Value = 96
$ws.Range("H30").Value = 96
$ws.Range("F32").Value = 91
$ws.Range("H32").Value = 91
$ws.Range("E33").Value = 244
$ws.Range("F33").Value = 123
$ws.Range("H33").Value = 123
$ws.Range("E34").Value = 174
$ws.Range("F34").Value = 105
$ws.Range("H34").Value = 105
$ws.Range("F35").Value = 73
$ws.Range("H35").Value = 73
$ws.Range("E36").Value = 52
$ws.Range("E37").Value = 129
$ws.Range("F37").Value = 62
$ws.Range("H37").Value = 62
$ws.Range("F38").Value = 52
$ws.Range("H38").Value = 52
$ws.Range("E39").Value = 160
$ws.Range("F39").Value = 76
$ws.Range("H39").Value = 76
$ws.Range("F40").Value = 98
$ws.Range("H40").Value = 98
$ws.Range("F41").Value = 149
$ws.Range("H41").Value = 149
$ws.Range("E42").Value = 292
$ws.Range("F42").Value = 156
$ws.Range("H42").Value = 156
$ws.Range("E43").Value = 98
$ws.Range("F43").Value = 50
$ws.Range("H43").Value = 50
$ws.Range("F44").Value = 121
$ws.Range("H44").Value = 121
$ws.Range("F45").Value = 55
$ws.Range("H45").Value = 55
$ws.Range("E46").Value = 254
$ws.Range("E47").Value = 365
$ws.Range("F47").Value = 179
$ws.Range("H47").Value = 179
$ws.Range("F48").Value = 68
$ws.Range("H48").Value = 68
$ws.Range("E49").Value = 246
$ws.Range("F49").Value = 103
$ws.Range("H49").Value = 103
$ws.Range("E50").Value = 209
$ws.Range("F50").Value = 89
$ws.Range("H50").Value = 89

$ws.Calculate()
